$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header - copy the header style/format from E1 (RXNO_DESC) onto F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "RXNO_DEF"

# New column data (RXNO_DEF) for rows 2-11
$defs = @(
  "['A generically dependent continuant that is about some thing. [IAO]']",
  "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]",
  "['B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]']",
  "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']",
  "['A processual entity that realizes a plan which is the concretization of a plan specification. [IAO]']",
  "['A directive information entity that describes an intended process endpoint. When part of a plan specification the concretization is realized in a planned process in which the bearer tries to effect the world so that the process endpoint is achieved. [IAO]']",
  "[]",
  "[]",
  "[]",
  "[]"
)

for ($i = 0; $i -lt $defs.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $defs[$i]
}
